# "27 - Analyse maken" use-case sheet update.
#
# Use case #14 ("Projecten aanmaken" / create projects) gets a second actor:
# besides "Studenten" (column D, already checked) the "Docenten" column (E)
# is now marked applicable too. Columns H:K hold helper formulas that build
# the human-readable use-case sentence from columns C:F, so this single
# input change ripples into the cached results of I19/J19/K19 automatically
# on recalculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19 = use case 14 "Projecten aanmaken": tick the "Docenten" checkbox.
$ws.Range("E19").Value = "x"

# Leave the selection/view the way the editing session left it: scrolled
# back to the top with H2:K24 (the generated use-case table) selected.
$ws.Range("H2:K24").Select() | Out-Null
